$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.228.89"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.905.29"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("D4").Value = "'1.002"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'307.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("E6").Value = "  +0.03%  "
$ws.Range("D7").Value = "'0.5254"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.41%  "
$ws.Range("D8").Value = "'0.3808"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.20%  "
$ws.Range("D9").Value = "'0.07294"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("D10").Value = "'21.71"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.47%  "
$ws.Range("D11").Value = "'0.9029"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.30%  "
$ws.Range("D12").Value = "'0.08217"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.15%  "
$ws.Range("D13").Value = "'96.32"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.56%  "
$ws.Range("D14").Value = "'5.352"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.13%  "
$ws.Range("D15").Value = "1.490.42"
$ws.Range("E15").Value = "  -21.87%  "
$ws.Range("D16").Value = "'1.002"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "'0.000008657"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.46%  "
$ws.Range("D18").Value = "'14.76"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.45%  "
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Value = "27.269.23"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("E21").Value = "  +0.50%  "
$ws.Range("D22").Value = "'10.81"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.61%  "
$ws.Range("D23").Value = "'6.503"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.98%  "
$ws.Range("B24").Value = "Monero"
$ws.Range("C24").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D24").Value = "'150.27"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.24%  "
$ws.Range("B25").Value = "LidoDAOToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D25").Value = "'2.338"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.35%  "
$ws.Range("D27").Value = "'1.743"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.95%  "
$ws.Range("D28").Value = "'116.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Value = "'4.841"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.43%  "
$ws.Range("E30").Value = "  -1.20%  "
$ws.Range("D31").Value = "'0.09250"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.73%  "
$ws.Range("D32").Value = "'0.8291"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.63%  "
$ws.Range("D33").Value = "'0.05055"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "'1.228"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.14%  "
$ws.Range("D35").Value = "'2.980"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.14%  "
$ws.Range("E36").Value = "  -2.38%  "
$ws.Range("D37").Value = "'2.725"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.84%  "
$ws.Range("D38").Value = "'0.5774"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.24%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").Value = "'9.212"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.20%  "
$ws.Range("D42").Value = "'6.614"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "'117.10"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.10%  "
$ws.Range("D44").Value = "'0.1520"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.41%  "
$ws.Range("D45").Value = "'0.4912"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.10%  "
$ws.Range("D46").Value = "'10.25"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.51%  "
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "'1.642"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.19%  "
$ws.Range("D49").Value = "'38.88"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.02%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").Value = "'0.06062"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.71%  "
